# po-cond01-test.xlsx edit: add the preliminary N-up-N-down staircase
# threshold-program columns (maskOnOff, stairUp, stairDn, dnDivUp, nRevs)
# to the single data sheet, and restore the view/selection state that
# Excel records after this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), O1:S1 -------------------------------------
# Shared strings must land in this exact order (maskOnOff, nRevs, stairUp,
# stairDn, dnDivUp) to reproduce the diff's sharedStrings.xml append order,
# so touch O1 then S1 before filling in P1:R1.
$ws.Range("O1").Value = "maskOnOff"
$ws.Range("S1").Value = "nRevs"
$ws.Range("P1").Value = "stairUp"
$ws.Range("Q1").Value = "stairDn"
$ws.Range("R1").Value = "dnDivUp"

# --- New data cells, rows 2-6: maskOnOff=1, stairUp=1, stairDn=1, --------
# --- dnDivUp=1, nRevs=4 for every condition row ---------------------------
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 15).Value = 1   # O: maskOnOff
    $ws.Cells.Item($r, 16).Value = 1   # P: stairUp
    $ws.Cells.Item($r, 17).Value = 1   # Q: stairDn
    $ws.Cells.Item($r, 18).Value = 1   # R: dnDivUp
    $ws.Cells.Item($r, 19).Value = 4   # S: nRevs
}

# --- View state: zoom + selection over the new staircase columns ---------
$excel.ActiveWindow.Zoom = 100
$ws.Range("P:R").Select()
